# daily auto push: 2026-01-12 22:32 UTC
# A new reading was logged for 2026/01/13 (Tue) and needs to be inserted
# into the timeline, which is kept in chronological order. Insert a new
# row at row 613 (shifting the existing rows 613:654 down to 614:655)
# and populate it with the new record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything from row 613 down by one row (xlShiftDown = -4121).
$ws.Rows("613:613").Insert(-4121)

# The date column stores plain text like "2026/01/13" (not a real Excel
# date), so force a text value instead of letting Excel auto-convert it
# to a date serial number, then clear the resulting quote-prefix style
# so the new cell matches the unstyled data rows around it.
$ws.Range("A613").Value = "'2026/01/13"
$ws.Range("A613").Style = "Normal"

$ws.Range("B613").Value = "火"
$ws.Range("C613").Value = 5
$ws.Range("D613").Value = 22
